$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume/change (E) columns for rows with new data
$ws.Range("D2").Value = "97.300.86"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "3.593.81"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "245.04"
$ws.Range("E5").Value = "  +3.34%  "
$ws.Range("D6").Value = "1.74"
$ws.Range("E6").Value = "  +16.23%  "
$ws.Range("D7").Value = "652.14"
$ws.Range("E7").Value = "  -0.87%  "
$ws.Range("D8").Value = "0.428"
$ws.Range("E8").Value = "  +6.29%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("E10").Value = "  +2.88%  "
$ws.Range("D11").Value = "3.591.66"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "44.50"
$ws.Range("E12").Value = "  +4.12%  "
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").Value = "6.48"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "4.260.58"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "97.227.84"
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("E17").Value = "  +2.78%  "
$ws.Range("D18").Value = "3.588.13"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").Value = "7.75"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "12.68"
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("D21").Value = "18.23"
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("D22").Value = "0.528"
$ws.Range("E22").Value = "  +7.76%  "
$ws.Range("D23").Value = "519.09"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("E25").Value = "  +4.29%  "
$ws.Range("E26").Value = "  -1.84%  "
$ws.Range("D27").Value = "104.06"
$ws.Range("E27").Value = "  +8.98%  "
$ws.Range("D28").Value = "13.25"
$ws.Range("E28").Value = "  +2.51%  "
$ws.Range("D31").Value = "2.98"
$ws.Range("E31").Value = "  -2.47%  "
$ws.Range("D32").Value = "11.99"
$ws.Range("E32").Value = "  +2.98%  "
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  +5.32%  "
$ws.Range("D35").Value = "0.990"
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("D36").Value = "31.92"
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("E37").Value = "  +3.42%  "
$ws.Range("D38").Value = "619.44"
$ws.Range("E38").Value = "  +3.25%  "
$ws.Range("D39").Value = "1.68"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("D40").Value = "8.77"
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("E41").Value = "  +1.78%  "
$ws.Range("E42").Value = "  +1.89%  "
$ws.Range("D43").Value = "0.932"
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("D45").Value = "6.12"
$ws.Range("E45").Value = "  +5.92%  "
$ws.Range("D46").Value = "0.446"
$ws.Range("E46").Value = "  +44.91%  "
$ws.Range("D47").Value = "0.0446"
$ws.Range("E47").Value = "  +6.82%  "
$ws.Range("D48").Value = "2.33"
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("D49").Value = "23.64"
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("D50").Value = "8.66"
$ws.Range("E50").Value = "  +5.09%  "
$ws.Range("D51").Value = "3.30"
$ws.Range("E51").Value = "  +6.76%  "

# Rows 29 and 30 swap: Hedera and WrappedeETH exchange positions with updated values
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "3.787.16"
$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "0.174"
$ws.Range("E30").Value = "  +19.92%  "
